$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: new diary entry (date + action) added alongside existing backlog cell H32 ---
$ws.Range("B31").Copy($ws.Range("B32"))
$ws.Range("C31").Copy($ws.Range("C32"))
$ws.Range("B32").Value = 42783
$ws.Range("C32").Value = "Pagination page links - SO question - (created Plunker)"

# --- Row 48: new backlog row (Pagination category, new item) ---
$ws.Range("G47").Copy($ws.Range("G48"))
$ws.Range("H47").Copy($ws.Range("H48"))
$ws.Range("H48").Value = "Go direct to certain page - eg pg 17"

# --- Row 33: new diary entry (date + action) added alongside existing backlog cell H33 ---
$ws.Range("B31").Copy($ws.Range("B33"))
$ws.Range("C31").Copy($ws.Range("C33"))
$ws.Range("B33").Value = 42783
$ws.Range("C33").Value = "Added wishlistController and sharedProperties to get & set carId"

# --- Row heights for the two newly-populated diary rows ---
$ws.Rows.Item(32).RowHeight = 15.75
$ws.Rows.Item(33).RowHeight = 15.75

# --- Selection moves to C37 ---
$ws.Range("C37").Select()
